$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("需求列表")
$iter = $wb.Worksheets.Item("迭代 1")

# Insert a new requirement row before the existing row 6 (pushes rows 6.. down by one).
$ws.Rows.Item(6).Insert() | Out-Null

# Fill in the new requirement row.
$ws.Cells.Item(6, 2).Value = "R-IL-INDOORQUICKLOCATION"
$ws.Cells.Item(6, 3).Value = "室内快捷定位"
$ws.Cells.Item(6, 4).Value = "1.附近卫生间定位。          2.停车场商场入口定位。               3.附近ATM定位。"
$ws.Rows.Item(6).RowHeight = 45

# New "对应场景" column header - copy the neighbouring header's formatting.
$ws.Cells.Item(1, 5).Value = "对应场景"
$ws.Range("D1").Copy() | Out-Null
$ws.Range("E1").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

# Move the selection / active sheet to match the edited workbook state.
$iter.Range("C18").Select() | Out-Null
$ws.Activate() | Out-Null
$ws.Range("E2").Select() | Out-Null
